$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lot code date in AC1
$ws.Range("AC1").Value = 45679

# Update the view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("AC2").Select()
